# Update the "want to go" count (column F) for a handful of events on both
# the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet, which
# mirrors the same rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    4  = 874
    5  = 40
    7  = 10611
    8  = 217
    19 = 299
    20 = 1014
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
